$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("L1").Value = "Date Deposited"
$ws.Range("L1").Style = $ws.Range("K1").Style

$ws.Range("L2").Select()
